$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the obsolete "path" column (A) values for the data rows (2-19).
# Using Clear() (not ClearContents) so the <c> elements are dropped entirely,
# matching the target workbook where those cells no longer exist.
$ws.Range("A2:A19").Clear()

# Update the price for the last article (Pex tube 20x2 isolé bleu 50m).
$ws.Range("D19").Value = 158.75

# Restore the cursor/selection to F19, as last left by the author.
$ws.Range("F19").Select()
